$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1.593749689512947
$ws.Range("D2").Value = 0.278418623906763
$ws.Range("E2").Value = 0.2587347637699864
$ws.Range("F2").Value = 0.8642576633317063
$ws.Range("G2").Value = 0.391028466203025
$ws.Range("H2").Value = 0.543973131779353
$ws.Range("I2").Value = 0.76626140285817
$ws.Range("J2").Value = 0.3294975375461036
$ws.Range("L2").Value = 0.7101507023506031
$ws.Range("O2").Value = 1.813793148630253
$ws.Range("B3").Value = 1.451906098937059
$ws.Range("D3").Value = 0.2812223158753611
$ws.Range("E3").Value = 0.2568917025592548
$ws.Range("F3").Value = 0.8739634813138508
$ws.Range("G3").Value = 0.3901157682615946
$ws.Range("H3").Value = 0.5478918876989312
$ws.Range("I3").Value = 0.7748952370540465
$ws.Range("J3").Value = 0.3207478604866765
$ws.Range("L3").Value = 0.6323783755559589
$ws.Range("O3").Value = 1.819763241671936
$ws.Range("B4").Value = 1.36451220110132
$ws.Range("D4").Value = 0.2830642500969525
$ws.Range("E4").Value = 0.255865125003929
$ws.Range("F4").Value = 0.8806136434229259
$ws.Range("G4").Value = 0.3899473382547853
$ws.Range("H4").Value = 0.5506297061354672
$ws.Range("I4").Value = 0.7808395279387206
$ws.Range("J4").Value = 0.315484411112152
$ws.Range("L4").Value = 0.5844217104742881
$ws.Range("O4").Value = 1.824943626498538
$ws.Range("B5").Value = 1.328825328967071
$ws.Range("D5").Value = 0.2838452123629196
$ws.Range("E5").Value = 0.2554733966294727
$ws.Range("F5").Value = 0.8834970669141455
$ws.Range("G5").Value = 0.3899769805831923
$ws.Range("H5").Value = 0.5518287550216385
$ws.Range("I5").Value = 0.7834235576553361
$ws.Range("J5").Value = 0.3133672823489491
$ws.Range("L5").Value = 0.564829254424609
$ws.Range("O5").Value = 1.827434870300181
$ws.Range("B6").Value = 1.322895211084244
$ws.Range("D6").Value = 0.2839767262131083
$ws.Range("E6").Value = 0.2554099632383071
$ws.Range("F6").Value = 0.8839863239938808
$ws.Range("G6").Value = 0.3899878300872572
$ws.Range("H6").Value = 0.5520328900899258
$ws.Range("I6").Value = 0.7838623984748274
$ws.Range("J6").Value = 0.3130174234542835
$ws.Range("L6").Value = 0.5615729863395131
$ws.Range("O6").Value = 1.827871481843488
$ws.Range("B7").Value = 1.36403120813435
$ws.Range("D7").Value = 0.2830746594122289
$ws.Range("E7").Value = 0.255859734021481
$ws.Range("F7").Value = 0.8806518284309632
$ws.Range("G7").Value = 0.3899473404737819
$ws.Range("H7").Value = 0.5506455394375607
$ws.Range("I7").Value = 0.7808737224485967
$ws.Range("J7").Value = 0.3154557458174025
$ws.Range("L7").Value = 0.5841576789528915
$ws.Range("O7").Value = 1.824975685815815
$ws.Range("B8").Value = 1.544906091053292
$ws.Range("D8").Value = 0.2793603780142142
$ws.Range("E8").Value = 0.2580775811597107
$ws.Range("F8").Value = 0.8674607043212603
$ws.Range("G8").Value = 0.3906322418552577
$ws.Range("H8").Value = 0.5452554531508014
$ws.Range("I8").Value = 0.7691048592430434
$ws.Range("J8").Value = 0.3264582778742664
$ws.Range("L8").Value = 0.6833781095581344
$ws.Range("O8").Value = 1.815536834910603
$ws.Range("B9").Value = 1.897112648831353
$ws.Range("D9").Value = 0.2730293410049889
$ws.Range("E9").Value = 0.2632526518531151
$ws.Range("F9").Value = 0.8470859134059623
$ws.Range("G9").Value = 0.3950990297723109
$ws.Range("H9").Value = 0.5373193501983309
$ws.Range("I9").Value = 0.75113071861945
$ws.Range("J9").Value = 0.3488825774652895
$ws.Range("L9").Value = 0.8762683784929948
$ws.Range("O9").Value = 1.809079370938662
$ws.Range("B10").Value = 2.154257551413821
$ws.Range("D10").Value = 0.2689544347052149
$ws.Range("E10").Value = 0.2675485544368712
$ws.Range("F10").Value = 0.8354822713029293
$ws.Range("G10").Value = 0.4003049261522023
$ws.Range("H10").Value = 0.5330975936867617
$ws.Range("I10").Value = 0.7410409988060138
$ws.Range("J10").Value = 0.3658560864531069
$ws.Range("L10").Value = 1.016892648135524
$ws.Range("O10").Value = 1.811731682370493
$ws.Range("B11").Value = 2.270867472782072
$ws.Range("D11").Value = 0.2672249375981792
$ws.Range("E11").Value = 0.2696081135827129
$ws.Range("F11").Value = 0.830937964694499
$ws.Range("G11").Value = 0.4030953295054758
$ws.Range("H11").Value = 0.531527085459615
$ws.Range("I11").Value = 0.7371285787447803
$ws.Range("J11").Value = 0.3736822037064087
$ws.Range("L11").Value = 1.080615899038321
$ws.Range("O11").Value = 1.814555104617881
$ws.Range("B12").Value = 2.314969724037837
$ws.Range("D12").Value = 0.2665878119952296
$ws.Range("E12").Value = 0.2704029865410575
$ws.Range("F12").Value = 0.8293230311652309
$ws.Range("G12").Value = 0.4042130169604548
$ws.Range("H12").Value = 0.5309827561584228
$ws.Range("I12").Value = 0.735744554064091
$ws.Range("J12").Value = 0.3766604675999616
$ws.Range("L12").Value = 1.104709272953357
$ws.Range("O12").Value = 1.815857563113354
$ws.Range("B13").Value = 2.305474021143652
$ws.Range("D13").Value = 0.266724237833607
$ws.Range("E13").Value = 0.2702311338679095
$ws.Range("F13").Value = 0.8296661212583629
$ws.Range("G13").Value = 0.403969584107486
$ws.Range("H13").Value = 0.5310977452411976
$ws.Range("I13").Value = 0.7360382893485422
$ws.Range("J13").Value = 0.3760183982191023
$ws.Range("L13").Value = 1.09952201616619
$ws.Range("O13").Value = 1.815566667315863
$ws.Range("B14").Value = 2.274496916961596
$ws.Range("D14").Value = 0.2671721645594545
$ws.Range("E14").Value = 0.2696732094325611
$ws.Range("F14").Value = 0.8308029790308282
$ws.Range("G14").Value = 0.4031860575917392
$ws.Range("H14").Value = 0.5314812929675696
$ws.Range("I14").Value = 0.7370127590342292
$ws.Range("J14").Value = 0.3739269351832206
$ws.Range("L14").Value = 1.082598831147266
$ws.Range("O14").Value = 1.814657578713962
$ws.Range("B15").Value = 2.255515253186445
$ws.Range("D15").Value = 0.2674488485248929
$ws.Range("E15").Value = 0.2693334077561715
$ws.Range("F15").Value = 0.8315131373611351
$ws.Range("G15").Value = 0.4027140809797629
$ws.Range("H15").Value = 0.5317227908072653
$ws.Range("I15").Value = 0.7376223536666373
$ws.Range("J15").Value = 0.3726477545295808
$ws.Range("L15").Value = 1.072227993289857
$ws.Range("O15").Value = 1.814131140099136
$ws.Range("B16").Value = 2.146629212177459
$ws.Range("D16").Value = 0.2690699552729114
$ws.Range("E16").Value = 0.2674160611709695
$ws.Range("F16").Value = 0.8357940541781801
$ws.Range("G16").Value = 0.4001310891891023
$ws.Range("H16").Value = 0.5332072779257686
$ws.Range("I16").Value = 0.7413103256696303
$ws.Range("J16").Value = 0.3653467068091487
$ws.Range("L16").Value = 1.012723062696239
$ws.Range("O16").Value = 1.81157976800344
$ws.Range("B17").Value = 2.079735342446099
$ws.Range("D17").Value = 0.2700962175434896
$ws.Range("E17").Value = 0.2662666775312061
$ws.Range("F17").Value = 0.8386085400067316
$ws.Range("G17").Value = 0.3986548669768553
$ws.Range("H17").Value = 0.5342076463820433
$ws.Range("I17").Value = 0.7437463616697784
$ws.Range("J17").Value = 0.3608943128599691
$ws.Range("L17").Value = 0.9761541559152533
$ws.Range("O17").Value = 1.810429264858982
$ws.Range("B18").Value = 2.04122541102015
$ws.Range("D18").Value = 0.2706981902673604
$ws.Range("E18").Value = 0.2656155076222362
$ws.Range("F18").Value = 0.8402964502243861
$ws.Range("G18").Value = 0.3978455056207935
$ws.Range("H18").Value = 0.5348159710275127
$ws.Range("I18").Value = 0.7452112624969942
$ws.Range("J18").Value = 0.3583432989690181
$ws.Range("L18").Value = 0.9550974961013594
$ws.Range("O18").Value = 1.809919667060939
$ws.Range("B19").Value = 2.02818079313181
$ws.Range("D19").Value = 0.270904018211894
$ws.Range("E19").Value = 0.2653967429752413
$ws.Range("F19").Value = 0.8408798047377815
$ws.Range("G19").Value = 0.3975782827334626
$ws.Range("H19").Value = 0.5350275947981231
$ws.Range("I19").Value = 0.7457181999437523
$ws.Range("J19").Value = 0.3574812791089812
$ws.Range("L19").Value = 0.9479641404785184
$ws.Range("O19").Value = 1.80977323235183
$ws.Range("B20").Value = 2.086859883377485
$ws.Range("D20").Value = 0.2699857603587716
$ws.Range("E20").Value = 0.2663880055034866
$ws.Range("F20").Value = 0.8383017801052262
$ws.Range("G20").Value = 0.3988078999678635
$ws.Range("H20").Value = 0.5340977459340905
$ws.Range("I20").Value = 0.743480442068055
$ws.Range("J20").Value = 0.3613672574997224
$ws.Range("L20").Value = 0.9800493908062435
$ws.Range("O20").Value = 1.810535986126098
$ws.Range("B21").Value = 2.283597167945175
$ws.Range("D21").Value = 0.2670401150979771
$ws.Range("E21").Value = 0.269836680738166
$ws.Range("F21").Value = 0.8304661798589024
$ws.Range("G21").Value = 0.4034145397526743
$ws.Range("H21").Value = 0.5313672677103085
$ws.Range("I21").Value = 0.7367238860008527
$ws.Range("J21").Value = 0.374540853119754
$ws.Range("L21").Value = 1.087570604562302
$ws.Range("O21").Value = 1.814918262431462
$ws.Range("B22").Value = 2.411852252924575
$ws.Range("D22").Value = 0.2652186740853182
$ws.Range("E22").Value = 0.2721777190160992
$ws.Range("F22").Value = 0.8259625526781633
$ws.Range("G22").Value = 0.4067810707419142
$ws.Range("H22").Value = 0.5298764700005734
$ws.Range("I22").Value = 0.7328765846874248
$ws.Range("J22").Value = 0.3832359687968818
$ws.Range("L22").Value = 1.157624124517099
$ws.Range("O22").Value = 1.819142484069658
$ws.Range("B23").Value = 2.343430631962804
$ws.Range("D23").Value = 0.2661813425170934
$ws.Range("E23").Value = 0.2709203506930393
$ws.Range("F23").Value = 0.8283096293688459
$ws.Range("G23").Value = 0.4049516305807259
$ws.Range("H23").Value = 0.5306452410574423
$ws.Range("I23").Value = 0.734877906692752
$ws.Range("J23").Value = 0.3785875328022144
$ws.Range("L23").Value = 1.120255714815073
$ws.Range("O23").Value = 1.816763225929066
$ws.Range("B24").Value = 2.083639039423133
$ws.Range("D24").Value = 0.2700356607907963
$ws.Range("E24").Value = 0.266333123128085
$ws.Range("F24").Value = 0.8384402487799392
$ws.Range("G24").Value = 0.3987385912394359
$ws.Range("H24").Value = 0.5341473285159282
$ws.Range("I24").Value = 0.7436004637530047
$ws.Range("J24").Value = 0.3611534120147724
$ws.Range("L24").Value = 0.9782884568734005
$ws.Range("O24").Value = 1.810487264521555
$ws.Range("B25").Value = 1.802108774896283
$ws.Range("D25").Value = 0.2746405076860832
$ws.Range("E25").Value = 0.2617652529089156
$ws.Range("F25").Value = 0.8520080017259986
$ws.Range("G25").Value = 0.3935541697566975
$ws.Range("H25").Value = 0.5391839687002147
$ws.Range("I25").Value = 0.7554463342407516
$ws.Range("J25").Value = 0.3488825774652895
$ws.Range("L25").Value = 0.8242734432087389
$ws.Range("O25").Value = 1.809530902347603